$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    'womens quick dry leggings',
    'womens race pants',
    'womens rainbow leggings',
    'womens rainbows',
    'womens recovery compression tights',
    'womens reebok clothing',
    'womens reflective pants',
    'womens reflective running gear',
    'womens riding tights',
    'womens rugby ball',
    'womens rugby shoes',
    'womens run leggings',
    'womens running baselayer',
    'womens running capris',
    'womens running clothes',
    'womens running clothing',
    'womens running compression',
    'womens running compression tights',
    'womens running equipment',
    'womens running gear',
    'womens running gear for cold weather',
    'womens running gear for winter',
    'womens running knee brace',
    'womens running leg warmers',
    'womens running leggings',
    'womens running leggings high waist',
    'womens running pants',
    'womens running shorts',
    'womens running tights',
    'womens running tights winter',
    'womens running tights with pockets',
    'womens running wear',
    'womens rustic clothing',
    'womens schooling tights',
    'womens sheer leggings',
    'womens shiny spandex shorts',
    'womens short tights',
    'womens shorts and capris',
    'womens shorts compression',
    'womens shorts knee length',
    'womens shorts mid rise',
    'womens sized leggings',
    'womens ski clothes',
    'womens ski pants',
    'womens ski pants black',
    'womens ski pants short length',
    'womens ski pants tight',
    'womens ski tight pants',
    'womens ski tights',
    'womens skii clothes',
    'womens skiing pants',
    'womens skirt with leggings',
    'womens skirts for leggings',
    'womens skirts with leggings',
    'womens snowboarding pants',
    'womens soccer',
    'womens soccer apparel',
    'womens soccer clothing',
    'womens soccer training pants',
    'womens soccer warm up pants',
    'womens soft leggings high waisted',
    'womens sport capri pants',
    'women''s sports compression pants & tights',
    'womens sports gear',
    'womens sports leggings',
    'womens sports tights',
    'womens stockings plus size',
    'womens stretch capri pants',
    'womens stretch capris',
    'womens stretch capris pants',
    'womens stretch yoga pants',
    'womens summer leggings pack',
    'womens support hose',
    'womens support panty hose',
    'womens support tights',
    'womens tall length pants',
    'womens tall tights',
    'womens tech suit',
    'womens tech suits',
    'womens tennis compression shorts',
    'womens tennis pants',
    'womens thick tights',
    'womens thick tights black',
    'womens thigh high',
    'womens thin tights',
    'womens tight',
    'womens tight pants',
    'womens tight running shorts',
    'womens tight shorts',
    'womens tight ski pants',
    'womens tight workout top',
    'womens tights',
    'womens tights 3 pack',
    'womens tights 4x',
    'womens tights and',
    'womens tights and leggings',
    'womens tights black',
    'women''s tights black',
    'womens tights gray',
    'womens tights small'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
